$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$shp.AlternativeText = "/*{{values:birthday.csv}}*/"
